# Auto-generated edit script: applies scheduled-runner price/profit updates
# to the per-class Leve profit tables (columns H-N) across all 8 sheets.
$wb = $excel.ActiveWorkbook

# ===== Sheet ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 112
$ws.Range("H112").Value = 1353.0513
$ws.Range("J112").Value = 1486.9395
$ws.Range("L112").Value = 4460.818499999999
$ws.Range("N112").Value = -6676.818499999999
# Row 113
$ws.Range("H113").Value = 2520.9
$ws.Range("I113").Value = 2449
$ws.Range("J113").Value = 2608.7778
$ws.Range("K113").Value = 2449
$ws.Range("L113").Value = 2608.7778
$ws.Range("M113").Value = 805
$ws.Range("N113").Value = -9116.7778
# Row 133
$ws.Range("H133").Value = 30184.445
$ws.Range("J133").Value = 30184.445
$ws.Range("L133").Value = 30184.445
$ws.Range("N133").Value = -40304.445
# Row 134
$ws.Range("H134").Value = 47737.5
$ws.Range("J134").Value = 47737.5
$ws.Range("L134").Value = 47737.5
$ws.Range("N134").Value = -57877.5
# Row 137
$ws.Range("H137").Value = 1119.3334
$ws.Range("I137").Value = 1244.2222
$ws.Range("J137").Value = 994.44446
$ws.Range("K137").Value = 3732.6666
$ws.Range("L137").Value = 2983.33338
$ws.Range("M137").Value = -1182.6666
$ws.Range("N137").Value = -8083.33338
# Row 138
$ws.Range("H138").Value = 4568805.5
$ws.Range("I138").Value = 7753644.5
$ws.Range("J138").Value = 3868.9666
$ws.Range("K138").Value = 23260933.5
$ws.Range("L138").Value = 11606.8998
$ws.Range("M138").Value = -23255793.5
$ws.Range("N138").Value = -21886.8998
# Row 140
$ws.Range("H140").Value = 28822.857
$ws.Range("J140").Value = 31626.666
$ws.Range("L140").Value = 31626.666
$ws.Range("N140").Value = -41986.666
# Row 141
$ws.Range("H141").Value = 1508.8334
$ws.Range("I141").Value = 1447.4375
$ws.Range("K141").Value = 4342.3125
$ws.Range("M141").Value = 837.6875

# ===== Sheet ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2852.9211
$ws.Range("I61").Value = 2952.1614
$ws.Range("J61").Value = 2413.4285
$ws.Range("K61").Value = 2952.1614
$ws.Range("L61").Value = 2413.4285
$ws.Range("M61").Value = -2740.1614
$ws.Range("N61").Value = -2837.4285
# Row 63
$ws.Range("H63").Value = 125003600
$ws.Range("I63").Value = 166670130
$ws.Range("J63").Value = 3998
$ws.Range("K63").Value = 166670130
$ws.Range("L63").Value = 3998
$ws.Range("M63").Value = -166669444
$ws.Range("N63").Value = -5370
# Row 66
$ws.Range("H66").Value = 125003600
$ws.Range("I66").Value = 166670130
$ws.Range("J66").Value = 3998
$ws.Range("K66").Value = 833350650
$ws.Range("L66").Value = 19990
$ws.Range("M66").Value = -833347218
$ws.Range("N66").Value = -26854
# Row 97
$ws.Range("H97").Value = 11852.333
$ws.Range("I97").Value = 14665.714
$ws.Range("K97").Value = 14665.714
$ws.Range("M97").Value = -14169.714
# Row 134
$ws.Range("H134").Value = 51221.35
$ws.Range("I134").Value = 11390
$ws.Range("J134").Value = 53710.812
$ws.Range("K134").Value = 11390
$ws.Range("L134").Value = 53710.812
$ws.Range("N134").Value = -63850.812
$ws.Range("M134").Value = -6320
# Row 136
$ws.Range("H136").Value = 2852.9211
$ws.Range("I136").Value = 2952.1614
$ws.Range("J136").Value = 2413.4285
$ws.Range("K136").Value = 8856.484199999999
$ws.Range("L136").Value = 7240.2855
$ws.Range("M136").Value = -6306.484199999999
$ws.Range("N136").Value = -12340.2855

# ===== Sheet BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 900
$ws.Range("I5").Value = 900
$ws.Range("K5").Value = 900
$ws.Range("M5").Value = -787
# Row 35
$ws.Range("H35").Value = 39800
$ws.Range("J35").Value = 39800
$ws.Range("L35").Value = 39800
$ws.Range("N35").Value = -40420
# Row 82
$ws.Range("H82").Value = 23818.176
$ws.Range("I82").Value = 14902.333
$ws.Range("J82").Value = 28681.363
$ws.Range("K82").Value = 14902.333
$ws.Range("L82").Value = 28681.363
$ws.Range("M82").Value = -14519.333
$ws.Range("N82").Value = -29447.363
# Row 85
$ws.Range("H85").Value = 23818.176
$ws.Range("I85").Value = 14902.333
$ws.Range("J85").Value = 28681.363
$ws.Range("K85").Value = 14902.333
$ws.Range("L85").Value = 28681.363
$ws.Range("M85").Value = -13576.333
$ws.Range("N85").Value = -31333.363
# Row 124
$ws.Range("H124").Value = 50780
$ws.Range("J124").Value = 50780
$ws.Range("L124").Value = 50780
$ws.Range("N124").Value = -60600

# ===== Sheet CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2719.7812
$ws.Range("I31").Value = 2070.7297
$ws.Range("K31").Value = 2070.7297
$ws.Range("M31").Value = -1775.7297
# Row 34
$ws.Range("H34").Value = 2719.7812
$ws.Range("I34").Value = 2070.7297
$ws.Range("K34").Value = 2070.7297
$ws.Range("M34").Value = -1868.7297
# Row 41
$ws.Range("H41").Value = 30009.285
$ws.Range("J41").Value = 30009.285
$ws.Range("L41").Value = 30009.285
$ws.Range("N41").Value = -30865.285
# Row 50
$ws.Range("H50").Value = 9244.75
$ws.Range("J50").Value = 9244.75
$ws.Range("L50").Value = 9244.75
$ws.Range("N50").Value = -10494.75
# Row 51
$ws.Range("H51").Value = 9265.4
$ws.Range("J51").Value = 9265.4
$ws.Range("L51").Value = 9265.4
$ws.Range("N51").Value = -10737.4
# Row 59
$ws.Range("H59").Value = 16020.25
$ws.Range("J59").Value = 16020.25
$ws.Range("L59").Value = 16020.25
$ws.Range("N59").Value = -18310.25
# Row 60
$ws.Range("H60").Value = 8000.6665
$ws.Range("J60").Value = 8251
$ws.Range("L60").Value = 8251
$ws.Range("N60").Value = -9273
# Row 61
$ws.Range("H61").Value = 9265.4
$ws.Range("J61").Value = 9265.4
$ws.Range("L61").Value = 9265.4
$ws.Range("N61").Value = -9961.4
# Row 68
$ws.Range("H68").Value = 17671.143
$ws.Range("J68").Value = 17671.143
$ws.Range("L68").Value = 17671.143
$ws.Range("N68").Value = -19169.143
# Row 71
$ws.Range("H71").Value = 17671.143
$ws.Range("J71").Value = 17671.143
$ws.Range("L71").Value = 53013.429
$ws.Range("N71").Value = -60501.429
# Row 74
$ws.Range("H74").Value = 13599.5
$ws.Range("J74").Value = 13599.5
$ws.Range("L74").Value = 13599.5
$ws.Range("N74").Value = -15347.5
# Row 77
$ws.Range("H77").Value = 13599.5
$ws.Range("J77").Value = 13599.5
$ws.Range("L77").Value = 40798.5
$ws.Range("N77").Value = -49534.5
# Row 99
$ws.Range("H99").Value = 2085.3809
$ws.Range("I99").Value = 1827.2727
$ws.Range("J99").Value = 2369.3
$ws.Range("K99").Value = 1827.2727
$ws.Range("L99").Value = 2369.3
$ws.Range("M99").Value = -329.2727
$ws.Range("N99").Value = -5365.3
# Row 126
$ws.Range("H126").Value = 2085.3809
$ws.Range("I126").Value = 1827.2727
$ws.Range("J126").Value = 2369.3
$ws.Range("K126").Value = 5481.8181
$ws.Range("L126").Value = 7107.900000000001
$ws.Range("M126").Value = -3011.8181
$ws.Range("N126").Value = -12047.9
# Row 135
$ws.Range("H135").Value = 52951.43
$ws.Range("J135").Value = 52951.43
$ws.Range("L135").Value = 52951.43
$ws.Range("N135").Value = -63091.43

# ===== Sheet CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 27778430
$ws.Range("I113").Value = 90909580
$ws.Range("J113").Value = 721.6
$ws.Range("K113").Value = 272728740
$ws.Range("L113").Value = 2164.8
$ws.Range("M113").Value = -272726570
$ws.Range("N113").Value = -6504.8

# ===== Sheet GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 3714.5264
$ws.Range("I126").Value = 2782.5715
$ws.Range("K126").Value = 8347.7145
$ws.Range("M126").Value = -5877.7145
# Row 141
$ws.Range("H141").Value = 64519.832
$ws.Range("J141").Value = 64519.832
$ws.Range("L141").Value = 64519.832
$ws.Range("N141").Value = -74879.83199999999

# ===== Sheet LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 6562.3667
$ws.Range("I7").Value = 4086.5293
$ws.Range("J7").Value = 9800
$ws.Range("K7").Value = 4086.5293
$ws.Range("L7").Value = 9800
$ws.Range("M7").Value = -3974.5293
$ws.Range("N7").Value = -10024
# Row 40
$ws.Range("H40").Value = 4555.5557
$ws.Range("I40").Value = 11418
$ws.Range("J40").Value = 1916.1538
$ws.Range("K40").Value = 11418
$ws.Range("L40").Value = 1916.1538
$ws.Range("M40").Value = -11282
$ws.Range("N40").Value = -2188.1538
# Row 126
$ws.Range("H126").Value = 6562.3667
$ws.Range("I126").Value = 4086.5293
$ws.Range("J126").Value = 9800
$ws.Range("K126").Value = 12259.5879
$ws.Range("L126").Value = 29400
$ws.Range("M126").Value = -9789.5879
$ws.Range("N126").Value = -34340
# Row 127
$ws.Range("H127").Value = 48195.75
$ws.Range("J127").Value = 48195.75
$ws.Range("L127").Value = 48195.75
$ws.Range("N127").Value = -58115.75
# Row 132
$ws.Range("H132").Value = 9809485
$ws.Range("I132").Value = 3736.5134
$ws.Range("J132").Value = 35724680
$ws.Range("K132").Value = 11209.5402
$ws.Range("L132").Value = 107174040
$ws.Range("M132").Value = -8679.540199999999
$ws.Range("N132").Value = -107179100
# Row 135
$ws.Range("H135").Value = 50154.758
$ws.Range("J135").Value = 50154.758
$ws.Range("L135").Value = 50154.758
$ws.Range("N135").Value = -60294.758

# ===== Sheet WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("N128").Value = 0
$ws.Range("L128").ClearContents()
# Row 135
$ws.Range("H135").Value = 42968.332
$ws.Range("J135").Value = 42968.332
$ws.Range("L135").Value = 42968.332
$ws.Range("N135").Value = -53108.332
# Row 137
$ws.Range("H137").Value = 70710
$ws.Range("J137").Value = 70710
$ws.Range("L137").Value = 70710
$ws.Range("N137").Value = -80910

